# Auto-generated script applying the scheduled-runner price update diff
# to the Belias_Profits workbook (columns H-N: market price / leve profit data).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 76924060
$ws.Range("I33").Value = 708.4545000000001
$ws.Range("J33").Value = 500002500
$ws.Range("K33").Value = 708.4545000000001
$ws.Range("L33").Value = 500002500
$ws.Range("M33").Value = -479.4545000000001
$ws.Range("N33").Value = -500002958
# Row 74
$ws.Range("H74").Value = 4090.9092
$ws.Range("I74").Value = 3785.7144
$ws.Range("J74").Value = 4625
$ws.Range("K74").Value = 3785.7144
$ws.Range("L74").Value = 4625
$ws.Range("M74").Value = -2849.7144
$ws.Range("N74").Value = -6497
# Row 76
$ws.Range("H76").Value = 3120.6897
$ws.Range("I76").Value = 3120.6897
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 3120.6897
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -2805.6897
$ws.Range("N76").Value = $null
# Row 77
$ws.Range("H77").Value = 4090.9092
$ws.Range("I77").Value = 3785.7144
$ws.Range("J77").Value = 4625
$ws.Range("K77").Value = 18928.572
$ws.Range("L77").Value = 23125
$ws.Range("M77").Value = -14248.572
$ws.Range("N77").Value = -32485
# Row 79
$ws.Range("H79").Value = 3120.6897
$ws.Range("I79").Value = 3120.6897
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 3120.6897
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -2028.6897
$ws.Range("N79").Value = $null
# Row 129
$ws.Range("H129").Value = 989.4761999999999
$ws.Range("I129").Value = 610.7273
$ws.Range("J129").Value = 1123.871
$ws.Range("K129").Value = 1832.1819
$ws.Range("L129").Value = 3371.613
$ws.Range("M129").Value = 3167.8181
$ws.Range("N129").Value = -13371.613

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 6
$ws.Range("H6").Value = 22222.445
$ws.Range("I6").Value = 28001
$ws.Range("J6").Value = 10665.333
$ws.Range("K6").Value = 28001
$ws.Range("L6").Value = 10665.333
$ws.Range("M6").Value = -27828
$ws.Range("N6").Value = -11011.333
# Row 45
$ws.Range("H45").Value = 1608.4166
$ws.Range("I45").Value = 1564.6
$ws.Range("J45").Value = 1639.7142
$ws.Range("K45").Value = 1564.6
$ws.Range("L45").Value = 1639.7142
$ws.Range("M45").Value = -1187.6
$ws.Range("N45").Value = -2393.7142
# Row 122
$ws.Range("H122").Value = 2117.238
$ws.Range("I122").Value = 1902.8334
$ws.Range("J122").Value = 2403.111
$ws.Range("K122").Value = 5708.5002
$ws.Range("L122").Value = 7209.333
$ws.Range("M122").Value = -3258.5002
$ws.Range("N122").Value = -12109.333

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Range("H7").Value = 35000
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").Value = $null
# Row 51
$ws.Range("H51").Value = 29900
$ws.Range("J51").Value = 29900
$ws.Range("L51").Value = 29900
$ws.Range("N51").Value = -30882
# Row 55
$ws.Range("H55").Value = 29126.334
$ws.Range("J55").Value = 29126.334
$ws.Range("L55").Value = 29126.334
$ws.Range("N55").Value = -29672.334

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 998.75
$ws.Range("I16").Value = 1060
$ws.Range("J16").Value = 962
$ws.Range("K16").Value = 1060
$ws.Range("L16").Value = 962
$ws.Range("M16").Value = -773
$ws.Range("N16").Value = -1536
# Row 28
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").Value = $null
# Row 113
$ws.Range("H113").Value = 998.75
$ws.Range("I113").Value = 1060
$ws.Range("J113").Value = 962
$ws.Range("K113").Value = 1060
$ws.Range("L113").Value = 962
$ws.Range("M113").Value = 1110
$ws.Range("N113").Value = -5302
# Row 122
$ws.Range("H122").Value = 1674
$ws.Range("I122").Value = 1565.3334
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 4696.0002
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -2246.0002
$ws.Range("N122").Value = -10900

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 173.83871
$ws.Range("I12").Value = 225.44444
$ws.Range("J12").Value = 152.72728
$ws.Range("K12").Value = 676.33332
$ws.Range("L12").Value = 458.18184
$ws.Range("M12").Value = -503.33332
$ws.Range("N12").Value = -804.18184
# Row 51
$ws.Range("H51").Value = 3269.2307
$ws.Range("I51").Value = 500
$ws.Range("K51").Value = 1500
$ws.Range("M51").Value = -1040
# Row 109
$ws.Range("H109").Value = 1951.5883
$ws.Range("I109").Value = 1028.1
$ws.Range("J109").Value = 3270.8572
$ws.Range("K109").Value = 3084.3
$ws.Range("L109").Value = 9812.571599999999
$ws.Range("M109").Value = -2044.3
$ws.Range("N109").Value = -11892.5716
# Row 113
$ws.Range("H113").Value = 685.8421
$ws.Range("I113").Value = 649.8333
$ws.Range("J113").Value = 747.5714
$ws.Range("K113").Value = 1949.4999
$ws.Range("L113").Value = 2242.7142
$ws.Range("M113").Value = 220.5001
$ws.Range("N113").Value = -6582.7142
# Row 115
$ws.Range("H115").Value = 2469.9
$ws.Range("J115").Value = 3450
$ws.Range("L115").Value = 10350
$ws.Range("N115").Value = -12700

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 3485.8
$ws.Range("J68").Value = 4333.3335
$ws.Range("L68").Value = 4333.3335
$ws.Range("N68").Value = -5831.3335
# Row 71
$ws.Range("H71").Value = 3485.8
$ws.Range("J71").Value = 4333.3335
$ws.Range("L71").Value = 21666.6675
$ws.Range("N71").Value = -29154.6675
# Row 96
$ws.Range("H96").Value = 27712.572
$ws.Range("J96").Value = 27712.572
$ws.Range("L96").Value = 27712.572
$ws.Range("N96").Value = -33204.572
# Row 98
$ws.Range("H98").Value = 26930
$ws.Range("J98").Value = 26930
$ws.Range("L98").Value = 26930
$ws.Range("N98").Value = -32920

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 93
$ws.Range("H93").Value = 30000
$ws.Range("J93").Value = 30000
$ws.Range("L93").Value = 30000
$ws.Range("N93").Value = -34992
# Row 98
$ws.Range("H98").Value = 30000
$ws.Range("J98").Value = 30000
$ws.Range("L98").Value = 30000
$ws.Range("N98").Value = -35990
# Row 126
$ws.Range("H126").Value = 1196.9722
$ws.Range("I126").Value = 667.1818
$ws.Range("J126").Value = 2029.5
$ws.Range("K126").Value = 2001.5454
$ws.Range("L126").Value = 6088.5
$ws.Range("M126").Value = 468.4546
$ws.Range("N126").Value = -11028.5
